# Primitives.docx edits:
#  - "datatypes" -> "data-types"   (Numbers paragraph)
#  - "Booleans exist true or false" -> "Booleans exist as true or false"
#  - "datatype" -> "data-type"     (String intro paragraph)
#  - "hello class" -> "hello class!"  (first, code-sample occurrence only)

$d = $word.ActiveDocument

# 1) "...int, double etc. datatypes but JavaScript..." -> "...data-types but JavaScript..."
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("datatypes but JavaScript", $true, $false, $false, $false, $false, $true, 1, $false, "data-types but JavaScript", 2)

# 2) "Booleans exist true or false" -> "Booleans exist as true or false"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Booleans exist true or false", $true, $false, $false, $false, $false, $true, 1, $false, "Booleans exist as true or false", 2)

# 3) "The next datatype we are going to focus" -> "The next data-type we are going to focus"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("The next datatype we are going to focus", $true, $false, $false, $false, $false, $true, 1, $false, "The next data-type we are going to focus", 2)

# 4) first curly-quoted "hello class" (the code line) -> "hello class!"
#    (a second, already-correct "hello class!" console-output line follows it,
#     so only replace the first match, not all)
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$openQuote = [char]0x201C
$closeQuote = [char]0x201D
$find.Execute($openQuote + "hello class" + $closeQuote, $true, $false, $false, $false, $false, $true, 1, $false, $openQuote + "hello class!" + $closeQuote, 1)
